$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data table (rows 2-13), columns A-T.
# A: Sending cluster, B: Ligand symbol, C: Receptor symbol, D: Target cluster
# E-J: ligand stats, K-P: receptor stats, Q-T: edge stats
$rows = @(
    @("ECs","Inhbb","Acvr1","ECs",3,1,2.525153666666667,7.575461000000001,0.7378191762484796,0.7378191762484798,3,1,7.818077,23.454231,0.1819018824829088,0.1819018824829088,19.74184580283233,177.676612225491,0.1342106970915875,0.1342106970915876),
    @("ECs","Inhbb","Acvr1","FAPs",3,1,2.525153666666667,7.575461000000001,0.7378191762484796,0.7378191762484798,3,1,22.011801,66.035403,0.5121448712693895,0.5121448712693895,55.58318000508701,500.2486200457831,0.3778703070398646,0.3778703070398647),
    @("ECs","Inhbb","Acvr1","M2",3,1,2.525153666666667,7.575461000000001,0.7378191762484796,0.7378191762484798,3,1,4.594576333333333,13.783729,0.1069012346955337,0.1069012346955337,11.60201127489656,104.418101474069,0.07887378092300405,0.07887378092300409),
    @("ECs","Inhbb","Acvr1","sCs",3,1,2.525153666666667,7.575461000000001,0.7378191762484796,0.7378191762484798,3,1,8.555183333333334,25.66555,0.1990520115521681,0.1990520115521681,21.60315256317222,194.42837306855,0.1468643911940235,0.1468643911940236),
    @("FAPs","Inhbb","Acvr1","ECs",2,0.6666666666666666,0.3739756666666667,1.121927,0.1092711394000877,0.1092711394000878,3,1,7.818077,23.454231,0.1819018824829088,0.1819018824829088,2.923770558126333,26.313935023137,0.0198766259579283,0.01987662595792831),
    @("FAPs","Inhbb","Acvr1","FAPs",2,0.6666666666666666,0.3739756666666667,1.121927,0.1092711394000877,0.1092711394000878,3,1,22.011801,66.035403,0.5121448712693895,0.5121448712693895,8.231877953509002,74.08690158158102,0.05596265362151744,0.05596265362151746),
    @("FAPs","Inhbb","Acvr1","M2",2,0.6666666666666666,0.3739756666666667,1.121927,0.1092711394000877,0.1092711394000878,3,1,4.594576333333333,13.783729,0.1069012346955337,0.1069012346955337,1.718259747309222,15.464337725783,0.01168121971845716,0.01168121971845716),
    @("FAPs","Inhbb","Acvr1","sCs",2,0.6666666666666666,0.3739756666666667,1.121927,0.1092711394000877,0.1092711394000878,3,1,8.555183333333334,25.66555,0.1990520115521681,0.1990520115521681,3.199430390538889,28.79487351485,0.02175064010218484,0.02175064010218485),
    @("sCs","Inhbb","Acvr1","ECs",3,1,0.5233266666666666,1.56998,0.1529096843514326,0.1529096843514326,3,1,7.818077,23.454231,0.1819018824829088,0.1819018824829088,4.091408176153332,36.82267358538,0.02781455943339297,0.02781455943339298),
    @("sCs","Inhbb","Acvr1","FAPs",3,1,0.5233266666666666,1.56998,0.1529096843514326,0.1529096843514326,3,1,22.011801,66.035403,0.5121448712693895,0.5121448712693895,11.51936244466,103.67426200194,0.07831191060800742,0.07831191060800743),
    @("sCs","Inhbb","Acvr1","M2",3,1,0.5233266666666666,1.56998,0.1529096843514326,0.1529096843514326,3,1,4.594576333333333,13.783729,0.1069012346955337,0.1069012346955337,2.404464317268888,21.64017885542,0.01634623405407247,0.01634623405407248),
    @("sCs","Inhbb","Acvr1","sCs",3,1,0.5233266666666666,1.56998,0.1529096843514326,0.1529096843514326,3,1,8.555183333333334,25.66555,0.1990520115521681,0.1990520115521681,4.477155576555555,40.294400189,0.03043698025595974,0.03043698025595975)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
